$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2,8).Value = "'8000.00"
$ws.Cells.Item(2,8).Style = "Normal"
$ws.Cells.Item(3,8).Value = "'5130.00"
$ws.Cells.Item(3,8).Style = "Normal"
$ws.Cells.Item(4,8).Value = "'23000.00"
$ws.Cells.Item(4,8).Style = "Normal"
$ws.Cells.Item(5,8).Value = "'31225.00"
$ws.Cells.Item(5,8).Style = "Normal"
$ws.Cells.Item(6,8).Value = "'578500.00"
$ws.Cells.Item(6,8).Style = "Normal"
$ws.Cells.Item(7,8).Value = "'2144000.00"
$ws.Cells.Item(7,8).Style = "Normal"
$ws.Cells.Item(8,8).Value = "'250000.00"
$ws.Cells.Item(8,8).Style = "Normal"
$ws.Cells.Item(9,8).Value = "'535.00"
$ws.Cells.Item(9,8).Style = "Normal"
$ws.Cells.Item(10,8).Value = "'13455.00"
$ws.Cells.Item(10,8).Style = "Normal"
$ws.Cells.Item(11,8).Value = "'850.00"
$ws.Cells.Item(11,8).Style = "Normal"
$ws.Cells.Item(12,8).Value = "'4225.00"
$ws.Cells.Item(12,8).Style = "Normal"
$ws.Cells.Item(13,8).Value = "'752275.54"
$ws.Cells.Item(13,8).Style = "Normal"
$ws.Cells.Item(14,8).Value = "'13000.00"
$ws.Cells.Item(14,8).Style = "Normal"
$ws.Cells.Item(15,8).Value = "'48400.00"
$ws.Cells.Item(15,8).Style = "Normal"
$ws.Cells.Item(16,8).Value = "'479500.00"
$ws.Cells.Item(16,8).Style = "Normal"
$ws.Cells.Item(17,8).Value = "'430.00"
$ws.Cells.Item(17,8).Style = "Normal"
$ws.Cells.Item(18,8).Value = "'1417860.54"
$ws.Cells.Item(18,8).Style = "Normal"
$ws.Cells.Item(19,8).Value = "'638231.93"
$ws.Cells.Item(19,8).Style = "Normal"
$ws.Cells.Item(20,8).Value = "'8973.56"
$ws.Cells.Item(20,8).Style = "Normal"
$ws.Cells.Item(21,8).Value = "'35750.50"
$ws.Cells.Item(21,8).Style = "Normal"
$ws.Cells.Item(22,8).Value = "'70120.85"
$ws.Cells.Item(22,8).Style = "Normal"
$ws.Cells.Item(23,8).Value = "'14847.00"
$ws.Cells.Item(23,8).Style = "Normal"
$ws.Cells.Item(24,8).Value = "'92150.00"
$ws.Cells.Item(24,8).Style = "Normal"
$ws.Cells.Item(25,8).Value = "'5769.80"
$ws.Cells.Item(25,8).Style = "Normal"
$ws.Cells.Item(26,8).Value = "'58656.76"
$ws.Cells.Item(26,8).Style = "Normal"
$ws.Cells.Item(27,8).Value = "'1470.00"
$ws.Cells.Item(27,8).Style = "Normal"
$ws.Cells.Item(28,8).Value = "'84223.94"
$ws.Cells.Item(28,8).Style = "Normal"
$ws.Cells.Item(29,8).Value = "'39900.00"
$ws.Cells.Item(29,8).Style = "Normal"
$ws.Cells.Item(30,8).Value = "'21600.00"
$ws.Cells.Item(30,8).Style = "Normal"
$ws.Cells.Item(31,8).Value = "'8400.00"
$ws.Cells.Item(31,8).Style = "Normal"
$ws.Cells.Item(32,8).Value = "'406.68"
$ws.Cells.Item(32,8).Style = "Normal"
$ws.Cells.Item(33,8).Value = "'75.92"
$ws.Cells.Item(33,8).Style = "Normal"
$ws.Cells.Item(34,8).Value = "'3581153.85"
$ws.Cells.Item(34,8).Style = "Normal"
$ws.Cells.Item(35,8).Value = "'9870.00"
$ws.Cells.Item(35,8).Style = "Normal"
$ws.Cells.Item(36,8).Value = "'37215.09"
$ws.Cells.Item(36,8).Style = "Normal"
$ws.Cells.Item(37,8).Value = "'11580.94"
$ws.Cells.Item(37,8).Style = "Normal"
$ws.Cells.Item(38,8).Value = "'26797.00"
$ws.Cells.Item(38,8).Style = "Normal"
$ws.Cells.Item(39,8).Value = "'181399.20"
$ws.Cells.Item(39,8).Style = "Normal"
$ws.Cells.Item(40,8).Value = "'11352.00"
$ws.Cells.Item(40,8).Style = "Normal"
$ws.Cells.Item(41,8).Value = "'10236.94"
$ws.Cells.Item(41,8).Style = "Normal"
$ws.Cells.Item(42,8).Value = "'51470.00"
$ws.Cells.Item(42,8).Style = "Normal"
$ws.Cells.Item(43,8).Value = "'1143.36"
$ws.Cells.Item(43,8).Style = "Normal"
$ws.Cells.Item(44,8).Value = "'43417.57"
$ws.Cells.Item(44,8).Style = "Normal"
$ws.Cells.Item(45,8).Value = "'210.00"
$ws.Cells.Item(45,8).Style = "Normal"
$ws.Cells.Item(46,8).Value = "'10365.00"
$ws.Cells.Item(46,8).Style = "Normal"
$ws.Cells.Item(47,8).Value = "'2269.77"
$ws.Cells.Item(47,8).Style = "Normal"
$ws.Cells.Item(48,8).Value = "'2676.46"
$ws.Cells.Item(48,8).Style = "Normal"
$ws.Cells.Item(49,8).Value = "'14405.00"
$ws.Cells.Item(49,8).Style = "Normal"
$ws.Cells.Item(50,8).Value = "'19036.74"
$ws.Cells.Item(50,8).Style = "Normal"
$ws.Cells.Item(51,8).Value = "'2767.40"
$ws.Cells.Item(51,8).Style = "Normal"
$ws.Cells.Item(52,8).Value = "'1066.04"
$ws.Cells.Item(52,8).Style = "Normal"
$ws.Cells.Item(53,8).Value = "'2715.00"
$ws.Cells.Item(53,8).Style = "Normal"
$ws.Cells.Item(54,8).Value = "'11250.00"
$ws.Cells.Item(54,8).Style = "Normal"
$ws.Cells.Item(55,8).Value = "'34396.92"
$ws.Cells.Item(55,8).Style = "Normal"
$ws.Cells.Item(56,8).Value = "'750.00"
$ws.Cells.Item(56,8).Style = "Normal"
$ws.Cells.Item(57,8).Value = "'79450.00"
$ws.Cells.Item(57,8).Style = "Normal"
$ws.Cells.Item(58,8).Value = "'127000.00"
$ws.Cells.Item(58,8).Style = "Normal"
$ws.Cells.Item(59,8).Value = "'5400.00"
$ws.Cells.Item(59,8).Style = "Normal"
$ws.Cells.Item(60,8).Value = "'56010.00"
$ws.Cells.Item(60,8).Style = "Normal"
$ws.Cells.Item(61,8).Value = "'48450.00"
$ws.Cells.Item(61,8).Style = "Normal"
$ws.Cells.Item(62,8).Value = "'38400.00"
$ws.Cells.Item(62,8).Style = "Normal"
$ws.Cells.Item(63,8).Value = "'40600.00"
$ws.Cells.Item(63,8).Style = "Normal"
$ws.Cells.Item(64,8).Value = "'7583.00"
$ws.Cells.Item(64,8).Style = "Normal"
$ws.Cells.Item(65,8).Value = "'1200.00"
$ws.Cells.Item(65,8).Style = "Normal"
$ws.Cells.Item(66,8).Value = "'550.00"
$ws.Cells.Item(66,8).Style = "Normal"
$ws.Cells.Item(67,8).Value = "'12100.00"
$ws.Cells.Item(67,8).Style = "Normal"
$ws.Cells.Item(68,8).Value = "'1035.00"
$ws.Cells.Item(68,8).Style = "Normal"
$ws.Cells.Item(69,8).Value = "'22250.00"
$ws.Cells.Item(69,8).Style = "Normal"
$ws.Cells.Item(70,8).Value = "'107000.00"
$ws.Cells.Item(70,8).Style = "Normal"
$ws.Cells.Item(71,8).Value = "'53006.00"
$ws.Cells.Item(71,8).Style = "Normal"
$ws.Cells.Item(72,8).Value = "'80203.00"
$ws.Cells.Item(72,8).Style = "Normal"
$ws.Cells.Item(73,8).Value = "'40750.00"
$ws.Cells.Item(73,8).Style = "Normal"
$ws.Cells.Item(74,8).Value = "'11557.14"
$ws.Cells.Item(74,8).Style = "Normal"
$ws.Cells.Item(75,8).Value = "'36500.00"
$ws.Cells.Item(75,8).Style = "Normal"
$ws.Cells.Item(76,8).Value = "'12070.00"
$ws.Cells.Item(76,8).Style = "Normal"
$ws.Cells.Item(77,8).Value = "'450.00"
$ws.Cells.Item(77,8).Style = "Normal"
$ws.Cells.Item(78,8).Value = "'350.00"
$ws.Cells.Item(78,8).Style = "Normal"
$ws.Cells.Item(79,8).Value = "'0.10"
$ws.Cells.Item(79,8).Style = "Normal"
$ws.Cells.Item(80,8).Value = "'4752.00"
$ws.Cells.Item(80,8).Style = "Normal"
$ws.Cells.Item(81,8).Value = "'9.60"
$ws.Cells.Item(81,8).Style = "Normal"
$ws.Cells.Item(82,8).Value = "'45000.00"
$ws.Cells.Item(82,8).Style = "Normal"
$ws.Cells.Item(83,8).Value = "'429.60"
$ws.Cells.Item(83,8).Style = "Normal"
$ws.Cells.Item(84,8).Value = "'21991.70"
$ws.Cells.Item(84,8).Style = "Normal"
$ws.Cells.Item(85,8).Value = "'40244.81"
$ws.Cells.Item(85,8).Style = "Normal"
$ws.Cells.Item(86,8).Value = "'1152.00"
$ws.Cells.Item(86,8).Style = "Normal"
$ws.Cells.Item(87,8).Value = "'2661.15"
$ws.Cells.Item(87,8).Style = "Normal"
$ws.Cells.Item(88,8).Value = "'1366.00"
$ws.Cells.Item(88,8).Style = "Normal"
$ws.Cells.Item(89,8).Value = "'900.00"
$ws.Cells.Item(89,8).Style = "Normal"
$ws.Cells.Item(90,8).Value = "'844.90"
$ws.Cells.Item(90,8).Style = "Normal"
$ws.Cells.Item(91,8).Value = "'86504.40"
$ws.Cells.Item(91,8).Style = "Normal"
$ws.Cells.Item(92,8).Value = "'13160.00"
$ws.Cells.Item(92,8).Style = "Normal"
$ws.Cells.Item(93,8).Value = "'2620.00"
$ws.Cells.Item(93,8).Style = "Normal"
$ws.Cells.Item(94,8).Value = "'2682.00"
$ws.Cells.Item(94,8).Style = "Normal"
$ws.Cells.Item(95,8).Value = "'2981.16"
$ws.Cells.Item(95,8).Style = "Normal"
$ws.Cells.Item(96,8).Value = "'2100.00"
$ws.Cells.Item(96,8).Style = "Normal"
$ws.Cells.Item(97,8).Value = "'10475.00"
$ws.Cells.Item(97,8).Style = "Normal"
$ws.Cells.Item(98,8).Value = "'12230.00"
$ws.Cells.Item(98,8).Style = "Normal"
$ws.Cells.Item(99,8).Value = "'13207.00"
$ws.Cells.Item(99,8).Style = "Normal"
$ws.Cells.Item(100,8).Value = "'32096.55"
$ws.Cells.Item(100,8).Style = "Normal"
$ws.Cells.Item(101,8).Value = "'110.00"
$ws.Cells.Item(101,8).Style = "Normal"
$ws.Cells.Item(102,8).Value = "'37285.00"
$ws.Cells.Item(102,8).Style = "Normal"
$ws.Cells.Item(103,8).Value = "'11339.33"
$ws.Cells.Item(103,8).Style = "Normal"
$ws.Cells.Item(104,8).Value = "'770.00"
$ws.Cells.Item(104,8).Style = "Normal"
$ws.Cells.Item(105,8).Value = "'215.00"
$ws.Cells.Item(105,8).Style = "Normal"
$ws.Cells.Item(106,8).Value = "'15535.00"
$ws.Cells.Item(106,8).Style = "Normal"
$ws.Cells.Item(107,8).Value = "'2850.00"
$ws.Cells.Item(107,8).Style = "Normal"
$ws.Cells.Item(108,8).Value = "'11900.00"
$ws.Cells.Item(108,8).Style = "Normal"
$ws.Cells.Item(109,8).Value = "'27100.00"
$ws.Cells.Item(109,8).Style = "Normal"
$ws.Cells.Item(110,8).Value = "'57000.00"
$ws.Cells.Item(110,8).Style = "Normal"
$ws.Cells.Item(111,8).Value = "'116000.00"
$ws.Cells.Item(111,8).Style = "Normal"
$ws.Cells.Item(112,8).Value = "'35000.00"
$ws.Cells.Item(112,8).Style = "Normal"
$ws.Cells.Item(113,8).Value = "'13500.00"
$ws.Cells.Item(113,8).Style = "Normal"
$ws.Cells.Item(114,8).Value = "'31000.00"
$ws.Cells.Item(114,8).Style = "Normal"
$ws.Cells.Item(115,8).Value = "'10000.00"
$ws.Cells.Item(115,8).Style = "Normal"
$ws.Cells.Item(116,8).Value = "'66000.00"
$ws.Cells.Item(116,8).Style = "Normal"
$ws.Cells.Item(117,8).Value = "'779000.00"
$ws.Cells.Item(117,8).Style = "Normal"
$ws.Cells.Item(118,8).Value = "'25100.00"
$ws.Cells.Item(118,8).Style = "Normal"
$ws.Cells.Item(119,8).Value = "'23709.62"
$ws.Cells.Item(119,8).Style = "Normal"
$ws.Cells.Item(120,8).Value = "'28333.52"
$ws.Cells.Item(120,8).Style = "Normal"
$ws.Cells.Item(121,8).Value = "'8780.00"
$ws.Cells.Item(121,8).Style = "Normal"
$ws.Cells.Item(122,8).Value = "'772.92"
$ws.Cells.Item(122,8).Style = "Normal"
$ws.Cells.Item(123,8).Value = "'197.22"
$ws.Cells.Item(123,8).Style = "Normal"
$ws.Cells.Item(124,8).Value = "'550.00"
$ws.Cells.Item(124,8).Style = "Normal"
$ws.Cells.Item(125,8).Value = "'5262.00"
$ws.Cells.Item(125,8).Style = "Normal"
$ws.Cells.Item(126,8).Value = "'240.00"
$ws.Cells.Item(126,8).Style = "Normal"
$ws.Cells.Item(127,8).Value = "'7000.00"
$ws.Cells.Item(127,8).Style = "Normal"
$ws.Cells.Item(128,8).Value = "'10000.00"
$ws.Cells.Item(128,8).Style = "Normal"
$ws.Cells.Item(129,8).Value = "'7000.00"
$ws.Cells.Item(129,8).Style = "Normal"
$ws.Cells.Item(130,8).Value = "'6050.00"
$ws.Cells.Item(130,8).Style = "Normal"
$ws.Cells.Item(131,8).Value = "'5000.00"
$ws.Cells.Item(131,8).Style = "Normal"
$ws.Cells.Item(132,8).Value = "'2000.00"
$ws.Cells.Item(132,8).Style = "Normal"
$ws.Cells.Item(133,8).Value = "'8000.00"
$ws.Cells.Item(133,8).Style = "Normal"
$ws.Cells.Item(134,8).Value = "'8500.00"
$ws.Cells.Item(134,8).Style = "Normal"
$ws.Cells.Item(135,8).Value = "'10000.00"
$ws.Cells.Item(135,8).Style = "Normal"
$ws.Cells.Item(136,8).Value = "'43504.40"
$ws.Cells.Item(136,8).Style = "Normal"
$ws.Cells.Item(137,8).Value = "'15000.00"
$ws.Cells.Item(137,8).Style = "Normal"
$ws.Cells.Item(138,8).Value = "'8000.00"
$ws.Cells.Item(138,8).Style = "Normal"
$ws.Cells.Item(139,8).Value = "'4860.00"
$ws.Cells.Item(139,8).Style = "Normal"
$ws.Cells.Item(140,8).Value = "'3800.00"
$ws.Cells.Item(140,8).Style = "Normal"
$ws.Cells.Item(141,8).Value = "'46113.00"
$ws.Cells.Item(141,8).Style = "Normal"
$ws.Cells.Item(142,8).Value = "'18500.00"
$ws.Cells.Item(142,8).Style = "Normal"
$ws.Cells.Item(143,8).Value = "'173700.00"
$ws.Cells.Item(143,8).Style = "Normal"
$ws.Cells.Item(144,8).Value = "'7260.00"
$ws.Cells.Item(144,8).Style = "Normal"
$ws.Cells.Item(145,8).Value = "'288.60"
$ws.Cells.Item(145,8).Style = "Normal"
$ws.Cells.Item(146,8).Value = "'11858.00"
$ws.Cells.Item(146,8).Style = "Normal"
$ws.Cells.Item(147,8).Value = "'12900.00"
$ws.Cells.Item(147,8).Style = "Normal"
$ws.Cells.Item(148,8).Value = "'9100.00"
$ws.Cells.Item(148,8).Style = "Normal"
$ws.Cells.Item(149,8).Value = "'2100.00"
$ws.Cells.Item(149,8).Style = "Normal"
$ws.Cells.Item(150,8).Value = "'31195.00"
$ws.Cells.Item(150,8).Style = "Normal"
$ws.Cells.Item(151,8).Value = "'30350.00"
$ws.Cells.Item(151,8).Style = "Normal"
$ws.Cells.Item(152,8).Value = "'123490.00"
$ws.Cells.Item(152,8).Style = "Normal"
$ws.Cells.Item(153,8).Value = "'3781.25"
$ws.Cells.Item(153,8).Style = "Normal"
$ws.Cells.Item(154,8).Value = "'6491.00"
$ws.Cells.Item(154,8).Style = "Normal"
$ws.Cells.Item(155,8).Value = "'9058.60"
$ws.Cells.Item(155,8).Style = "Normal"
$ws.Cells.Item(156,8).Value = "'360.00"
$ws.Cells.Item(156,8).Style = "Normal"
$ws.Cells.Item(157,8).Value = "'297.32"
$ws.Cells.Item(157,8).Style = "Normal"
$ws.Cells.Item(158,8).Value = "'166145.10"
$ws.Cells.Item(158,8).Style = "Normal"
$ws.Cells.Item(159,8).Value = "'3042.00"
$ws.Cells.Item(159,8).Style = "Normal"
$ws.Cells.Item(160,8).Value = "'10925.00"
$ws.Cells.Item(160,8).Style = "Normal"
$ws.Cells.Item(161,8).Value = "'274145.00"
$ws.Cells.Item(161,8).Style = "Normal"
$ws.Cells.Item(162,8).Value = "'32675.38"
$ws.Cells.Item(162,8).Style = "Normal"
$ws.Cells.Item(163,8).Value = "'11985.00"
$ws.Cells.Item(163,8).Style = "Normal"
$ws.Cells.Item(164,8).Value = "'10704.00"
$ws.Cells.Item(164,8).Style = "Normal"
$ws.Cells.Item(165,8).Value = "'11160.00"
$ws.Cells.Item(165,8).Style = "Normal"
$ws.Cells.Item(166,8).Value = "'6800.00"
$ws.Cells.Item(166,8).Style = "Normal"
$ws.Cells.Item(167,8).Value = "'92641.11"
$ws.Cells.Item(167,8).Style = "Normal"
$ws.Cells.Item(168,8).Value = "'60000.00"
$ws.Cells.Item(168,8).Style = "Normal"
$ws.Cells.Item(169,8).Value = "'30000.00"
$ws.Cells.Item(169,8).Style = "Normal"
$ws.Cells.Item(170,8).Value = "'30000.00"
$ws.Cells.Item(170,8).Style = "Normal"
$ws.Cells.Item(171,8).Value = "'75000.00"
$ws.Cells.Item(171,8).Style = "Normal"
$ws.Cells.Item(172,8).Value = "'30000.00"
$ws.Cells.Item(172,8).Style = "Normal"
$ws.Cells.Item(173,8).Value = "'30000.00"
$ws.Cells.Item(173,8).Style = "Normal"
$ws.Cells.Item(174,8).Value = "'30000.00"
$ws.Cells.Item(174,8).Style = "Normal"
$ws.Cells.Item(175,8).Value = "'30000.00"
$ws.Cells.Item(175,8).Style = "Normal"
$ws.Cells.Item(176,8).Value = "'60000.00"
$ws.Cells.Item(176,8).Style = "Normal"
$ws.Cells.Item(177,8).Value = "'60000.00"
$ws.Cells.Item(177,8).Style = "Normal"
$ws.Cells.Item(178,8).Value = "'13200.00"
$ws.Cells.Item(178,8).Style = "Normal"
$ws.Cells.Item(179,8).Value = "'480000.00"
$ws.Cells.Item(179,8).Style = "Normal"
$ws.Cells.Item(180,8).Value = "'21888.00"
$ws.Cells.Item(180,8).Style = "Normal"
$ws.Cells.Item(181,8).Value = "'10052200.19"
$ws.Cells.Item(181,8).Style = "Normal"
$ws.Cells.Item(182,8).Value = "'6395.00"
$ws.Cells.Item(182,8).Style = "Normal"
$ws.Cells.Item(183,8).Value = "'6047027.00"
$ws.Cells.Item(183,8).Style = "Normal"
$ws.Cells.Item(184,8).Value = "'201500.00"
$ws.Cells.Item(184,8).Style = "Normal"
$ws.Cells.Item(185,8).Value = "'266500.00"
$ws.Cells.Item(185,8).Style = "Normal"
$ws.Cells.Item(186,8).Value = "'201500.00"
$ws.Cells.Item(186,8).Style = "Normal"
$ws.Cells.Item(187,8).Value = "'201500.00"
$ws.Cells.Item(187,8).Style = "Normal"
$ws.Cells.Item(188,8).Value = "'201500.00"
$ws.Cells.Item(188,8).Style = "Normal"
$ws.Cells.Item(189,8).Value = "'1500.00"
$ws.Cells.Item(189,8).Style = "Normal"
$ws.Cells.Item(190,8).Value = "'208000.00"
$ws.Cells.Item(190,8).Style = "Normal"
$ws.Cells.Item(191,8).Value = "'377700.00"
$ws.Cells.Item(191,8).Style = "Normal"
$ws.Cells.Item(192,8).Value = "'201500.00"
$ws.Cells.Item(192,8).Style = "Normal"
$ws.Cells.Item(193,8).Value = "'305400.00"
$ws.Cells.Item(193,8).Style = "Normal"
$ws.Cells.Item(194,8).Value = "'380800.00"
$ws.Cells.Item(194,8).Style = "Normal"
$ws.Cells.Item(195,8).Value = "'171200.00"
$ws.Cells.Item(195,8).Style = "Normal"
$ws.Cells.Item(196,8).Value = "'338300.00"
$ws.Cells.Item(196,8).Style = "Normal"
$ws.Cells.Item(197,8).Value = "'201500.00"
$ws.Cells.Item(197,8).Style = "Normal"
$ws.Cells.Item(198,8).Value = "'396500.00"
$ws.Cells.Item(198,8).Style = "Normal"
$ws.Cells.Item(199,8).Value = "'403000.00"
$ws.Cells.Item(199,8).Style = "Normal"
$ws.Cells.Item(200,8).Value = "'256300.00"
$ws.Cells.Item(200,8).Style = "Normal"
$ws.Cells.Item(201,8).Value = "'416900.00"
$ws.Cells.Item(201,8).Style = "Normal"
$ws.Cells.Item(202,8).Value = "'743300.00"
$ws.Cells.Item(202,8).Style = "Normal"
$ws.Cells.Item(203,8).Value = "'367700.00"
$ws.Cells.Item(203,8).Style = "Normal"
$ws.Cells.Item(204,8).Value = "'21000.00"
$ws.Cells.Item(204,8).Style = "Normal"
$ws.Cells.Item(205,8).Value = "'533900.00"
$ws.Cells.Item(205,8).Style = "Normal"
$ws.Cells.Item(206,8).Value = "'403000.00"
$ws.Cells.Item(206,8).Style = "Normal"
$ws.Cells.Item(207,8).Value = "'212470.00"
$ws.Cells.Item(207,8).Style = "Normal"
$ws.Cells.Item(208,8).Value = "'922020.00"
$ws.Cells.Item(208,8).Style = "Normal"
$ws.Cells.Item(209,8).Value = "'20400.00"
$ws.Cells.Item(209,8).Style = "Normal"
$ws.Cells.Item(210,8).Value = "'628579.16"
$ws.Cells.Item(210,8).Style = "Normal"
$ws.Cells.Item(211,8).Value = "'74764.20"
$ws.Cells.Item(211,8).Style = "Normal"
$ws.Cells.Item(212,8).Value = "'1400.00"
$ws.Cells.Item(212,8).Style = "Normal"
$ws.Cells.Item(213,8).Value = "'70500.00"
$ws.Cells.Item(213,8).Style = "Normal"
$ws.Cells.Item(214,8).Value = "'208752.00"
$ws.Cells.Item(214,8).Style = "Normal"
$ws.Cells.Item(215,8).Value = "'11000.00"
$ws.Cells.Item(215,8).Style = "Normal"
$ws.Cells.Item(216,8).Value = "'118000.00"
$ws.Cells.Item(216,8).Style = "Normal"
$ws.Cells.Item(217,8).Value = "'7000.00"
$ws.Cells.Item(217,8).Style = "Normal"
$ws.Cells.Item(218,8).Value = "'5000.00"
$ws.Cells.Item(218,8).Style = "Normal"
$ws.Cells.Item(219,8).Value = "'9600.00"
$ws.Cells.Item(219,8).Style = "Normal"
$ws.Cells.Item(220,8).Value = "'96882.04"
$ws.Cells.Item(220,8).Style = "Normal"
$ws.Cells.Item(221,8).Value = "'15500.00"
$ws.Cells.Item(221,8).Style = "Normal"
$ws.Cells.Item(222,8).Value = "'70500.00"
$ws.Cells.Item(222,8).Style = "Normal"
$ws.Cells.Item(223,8).Value = "'380.00"
$ws.Cells.Item(223,8).Style = "Normal"
$ws.Cells.Item(224,8).Value = "'48985.00"
$ws.Cells.Item(224,8).Style = "Normal"